$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Trends Status sheet: update values
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("E2").Value = 5.3
$ws1.Range("C3").Value = 7
$ws1.Range("E3").Value = 36.8
$ws1.Range("C4").Value = 10
$ws1.Range("E4").Value = 52.6
$ws1.Range("C5").Value = 1
$ws1.Range("E5").Value = 5.3
$ws1.Range("C7").Value = 19
$ws1.Range("B8").Value = 386
$ws1.Range("C8").Value = 348

# ---------------------------------------------------------------------------
# 2. Priority Status sheet: update values
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Range("B2").Value = 103
$ws3.Range("B3").Value = 286
$ws3.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# 3. Species qualification sheet: update text + values
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("A2").Value = "SoIB Assessment"
$ws4.Range("B2").Value = 386
$ws4.Range("B4").Value = 38
$ws4.Range("C4").Value = 19

# ---------------------------------------------------------------------------
# 4. High Priority break-up sheet: rename to "Interannual update - High Pri"
#    and replace its contents with the new break-up table, then add a new
#    sheet "Major update - High Priority " that keeps the old single-row
#    "IUCN" content.
# ---------------------------------------------------------------------------
$oldSheet = $wb.Worksheets.Item("High Priority break-up")

# Create the new sheet that will carry the old content forward, placed right
# after the existing break-up sheet.
$newSheet = $wb.Worksheets.Add($null, $oldSheet)
$newSheet.Name = "Major update - High Priority "

# Reuse the existing header formatting (bold + centered) for the new sheet.
$oldSheet.Range("A1:E1").Copy()
$newSheet.Range("A1:E1").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "Break-up"
$newSheet.Range("B1").Value = "High Species (no.)"
$newSheet.Range("C1").Value = "High Species (perc.)"
$newSheet.Range("D1").Value = "New High Species (no.)"
$newSheet.Range("E1").Value = "New High Species (perc.)"

$newSheet.Range("A2").Value = "IUCN"
$newSheet.Range("B2").Value = 10
$newSheet.Range("C2").Value = 100
$newSheet.Range("D2").Value = 10
$newSheet.Range("E2").Value = 100

# Rename the original sheet and overwrite it with the new interannual data.
$oldSheet.Name = "Interannual update - High Pri"

$oldSheet.Range("A1").Value = "Break-up"
$oldSheet.Range("B1").Value = "High Species (no.)"
$oldSheet.Range("C1").Value = "High Species (perc.)"
$oldSheet.Range("D1").Value = "New High Species (no.)"
$oldSheet.Range("E1").Value = "New High Species (perc.)"

$oldSheet.Range("A2").Value = "Trend New"
$oldSheet.Range("B2").Value = 73
$oldSheet.Range("C2").Value = 70.90000000000001
$oldSheet.Range("D2").Value = 73
$oldSheet.Range("E2").Value = 77.7

$oldSheet.Range("A3").Value = "IUCN"
$oldSheet.Range("B3").Value = 30
$oldSheet.Range("C3").Value = 29.1
$oldSheet.Range("D3").Value = 21
$oldSheet.Range("E3").Value = 22.3
